$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, pushing existing rows 150-185 down to 151-186.
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new weekly record
# (same market/category as its neighbours, new date + origin).
$ws.Range("A150").Value = 7
$ws.Range("B150").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C150").Value = "Ñuble"
$ws.Range("D150").Value = 44511
$ws.Range("E150").Value = 16
$ws.Range("F150").Value = 100112008
$ws.Range("G150").Value = "Coliflor"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 300
$ws.Range("K150").Value = 700
$ws.Range("L150").Value = 750
$ws.Range("M150").Value = 725
$ws.Range("N150").Value = "$/unidad"
$ws.Range("O150").Value = "Región del Maule"
$ws.Range("P150").Value = 725
$ws.Range("Q150").Value = 1
$ws.Range("R150").Value = "Hortaliza"
